# Updates cryptos list values (Price / Volume(1h) columns) to match the
# latest scraped snapshot, per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.842.09'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '1.560.24'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''205.60'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '''0.247'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").Value = '''21.52'
$ws.Range("E9").Value = '  -3.42%  '
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("D11").Value = '''0.0862'
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("D12").Value = '1.781.23'
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("D13").Value = '1.558.62'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("E15").Value = '  -0.94%  '
$ws.Range("D16").Value = '26.837.47'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").Value = '''61.31'
$ws.Range("E17").Value = '  -2.64%  '
$ws.Range("D18").Value = '''214.70'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '''9.12'
$ws.Range("E23").Value = '  -2.37%  '
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("D25").Value = '''153.56'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").Value = '''14.97'
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '''0.102'
$ws.Range("E29").Value = '  -1.22%  '
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").Value = '1.375.97'
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("D34").Value = '''2.92'
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("E35").Value = '  -2.58%  '
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("D37").Value = '''0.921'
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("E39").Value = '  +1.59%  '
$ws.Range("D40").Value = '''0.809'
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").Value = '''5.55'
$ws.Range("E43").Value = '  +4.63%  '
$ws.Range("D44").Value = '''1.77'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").Value = '''63.46'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '1.695.40'
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").Value = '''86.51'
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("D49").Value = '''0.0512'
$ws.Range("E49").Value = '  +3.92%  '
$ws.Range("D50").Value = '0.0₇0977'
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("E51").Value = '  +0.30%  '
